# Add season record columns (Wins, Losses, Ties) to the roster/stats sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new headers Wins / Losses / Ties in AC1:AE1 ---
$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"

# Copy the formatting of an existing header cell (AB1) onto the new header
# cells so they keep the same bold/border/alignment style used by the rest
# of the header row.
$ws.Range("AB1").Copy()
$ws.Range("AC1:AE1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Re-assert the header text (PasteSpecial only copies formats, but make sure
# values are correct regardless of paste semantics).
$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"

# --- Data rows (2-39): every team/season row gets the same record ---
$wins = 75
$losses = 87
$ties = 0

for ($r = 2; $r -le 39; $r++) {
    $ws.Cells.Item($r, 29).Value = $wins    # column AC
    $ws.Cells.Item($r, 30).Value = $losses  # column AD
    $ws.Cells.Item($r, 31).Value = $ties    # column AE
}
